$wb = $excel.ActiveWorkbook

# xlLineStyleNone constant (used to explicitly clear an edge we don't want)
$xlNone = -4142
$xlContinuous = 1

function Set-TopBottomBorder($range) {
    # top + bottom thin, left/right none  -> reuses existing borderId 4
    $range.Borders.Item(8).LineStyle = $xlContinuous   # xlEdgeTop
    $range.Borders.Item(9).LineStyle = $xlContinuous   # xlEdgeBottom
    $range.Borders.Item(7).LineStyle = $xlNone         # xlEdgeLeft
    $range.Borders.Item(10).LineStyle = $xlNone        # xlEdgeRight
}

function Set-RightTopBottomBorder($range) {
    # right + top + bottom thin, left none -> reuses existing borderId 5
    $range.Borders.Item(8).LineStyle = $xlContinuous   # xlEdgeTop
    $range.Borders.Item(9).LineStyle = $xlContinuous   # xlEdgeBottom
    $range.Borders.Item(10).LineStyle = $xlContinuous  # xlEdgeRight
    $range.Borders.Item(7).LineStyle = $xlNone         # xlEdgeLeft
}

# --- Sheet 1: quality_comparison ---
$ws1 = $wb.Worksheets.Item("quality_comparison")

Set-TopBottomBorder      $ws1.Range("C1")
Set-RightTopBottomBorder $ws1.Range("D1")

$ws1.Range("C2").Value = "approach"

# --- Sheet 2: computational_comparison ---
$ws2 = $wb.Worksheets.Item("computational_comparison")

Set-TopBottomBorder      $ws2.Range("C1")
Set-RightTopBottomBorder $ws2.Range("D1")
Set-TopBottomBorder      $ws2.Range("F1")
Set-RightTopBottomBorder $ws2.Range("G1")

$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# G5 had an empty inline-string placeholder cell; it should be removed entirely.
$ws2.Range("G5").ClearContents()
